$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values for the refreshed symbol list.
# Values are forced to Text format to preserve exact formatting (trailing zeros,
# percent signs, thousands separators) exactly as scraped, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.22%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.96"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "9.60%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.275"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.22%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07500"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.89%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.868"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.54%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.815"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "7.22%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.478"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "6.16%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9212"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.89%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1695"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.71%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07866"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.26%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08021"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.00%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03061"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "5.12%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09911"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "9.76%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001496"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-6.67%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04609"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.87%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006146"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.53%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.459"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.98%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.228"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.22%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.85%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.39%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.503"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "12.35%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001214"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.38%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.72%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001398"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "19.71%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "16.07%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01724"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2,548.01%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04497"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.73%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006887"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.43%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1349"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "7.76%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002197"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.33%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01283"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "9.63%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006173"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.37%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.866"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.31%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01498"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "15.45%"
